# Updates the cryptos list (Price / Volume(1h) columns) with refreshed
# market data, and fixes the ranking order of two swapped pairs of rows
# (Polkadot/WrappedBTC at rows 16-17, CoreDAO/USDe at rows 49-50).
#
# Column D ("Price") values are free-form text (e.g. "66.438.53",
# "1.00"), not real numbers, so for any value that looks like a plain
# number we force the cell to Text format before writing, then restore
# the "Normal" style afterwards so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.438.53'
$ws.Range('E2').Value = '  +2.96%  '

# Row 3
$ws.Range('D3').Value = '3.188.72'
$ws.Range('E3').Value = '  +1.49%  '

# Row 4
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.57%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.34'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.90%  '

# Row 7
$ws.Range('E7').Value = '  +0.00%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.569'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +8.32%  '

# Row 9
$ws.Range('D9').Value = '3.185.65'
$ws.Range('E9').Value = '  +1.38%  '

# Row 10
$ws.Range('E10').Value = '  +1.94%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.91'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.38%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.520'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.95%  '

# Row 13
$ws.Range('E13').Value = '  +2.77%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '39.24'
$ws.Range('D14').Style = 'Normal'

# Row 15
$ws.Range('D15').Value = '3.710.70'
$ws.Range('E15').Value = '  +1.46%  '

# Row 16
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.50'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.22%  '

# Row 17
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '66.446.02'
$ws.Range('E17').Value = '  +2.78%  '

# Row 18
$ws.Range('D18').Value = '3.190.41'
$ws.Range('E18').Value = '  +1.45%  '

# Row 19
$ws.Range('E19').Value = '  +0.66%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '520.09'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.46%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.44'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.76%  '

# Row 23
$ws.Range('E23').Value = '  +5.63%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.96'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.41%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.27'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.47%  '

# Row 26
$ws.Range('E26').Value = '  +0.14%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.27'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.75%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.80%  '

# Row 29
$ws.Range('E29').Value = '  +8.29%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.13'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +15.60%  '

# Row 31
$ws.Range('E31').Value = '  +4.23%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.40'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.24%  '

# Row 33
$ws.Range('E33').Value = '  +3.07%  '

# Row 34
$ws.Range('E34').Value = '  +0.07%  '

# Row 35
$ws.Range('E35').Value = '  +1.60%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '510.94'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.48%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '54.95'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.84%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0904'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.90%  '

# Row 39
$ws.Range('E39').Value = '  +3.13%  '

# Row 40
$ws.Range('E40').Value = '  +10.71%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.93'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.89%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.88'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.29%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.302'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.78%  '

# Row 44
$ws.Range('D44').Value = '0.0₃0669'
$ws.Range('E44').Value = '  +15.80%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.46'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.08%  '

# Row 46
$ws.Range('D46').Value = '2.906.64'
$ws.Range('E46').Value = '  -2.88%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.76'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.73%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.119'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.51%  '

# Row 49
$ws.Range('B49').Value = 'CoreDAO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.68'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +9.69%  '

# Row 50
$ws.Range('B50').Value = 'USDe'
$ws.Range('C50').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.999'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.01%  '

# Row 51
$ws.Range('E51').Value = '  +6.20%  '
